$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper idea: Word/this engine will coalesce adjacent runs that share an
# identical rPr whenever a Range.Text assignment touches that neighbourhood.
# To produce an explicit run split at an exact boundary (matching the target
# OOXML) we first let a plain text replace merge things, then re-impose a
# split by toggling Font.Bold off/on (or on/off, matching the original
# boolean) on the sub-range we want isolated -- this forces the engine to
# keep that sub-range as its own run while leaving rPr unchanged (no stray
# w:val="0" residue) because we return the value to what it already was.
# ---------------------------------------------------------------------------

# =====================================================================
# Change 1: "Double " + "First Class" + " Undergraduate Degree..." (bold)
#   -> "Double First" | "-" | "Class Undergraduate Degree..."
#   (also drops the gramStart/gramEnd proofErr markers around "First Class")
# =====================================================================
$r1 = $d.Content
$r1.Find.Execute("Double First Class Undergraduate Degree")
$r1.Text = "Double First-Class Undergraduate Degree"

$r1b = $d.Content
$r1b.Find.Execute("Double First-Class Undergraduate")
$start1 = $r1b.Start
$hyphen1 = $d.Range($start1 + 12, $start1 + 13)
$hyphen1.Font.Bold = $false
$hyphen1.Font.Bold = $true

# =====================================================================
# Change 2: " as a high first. Achieved official " + "first class" +
#           " results in second and third year." (normal weight)
#   -> " as a high first. Achieved official first" | "-" |
#      "class results in second " | "year " | "and third year" | "."
#   (also drops the gramStart/gramEnd proofErr markers around "first class")
# =====================================================================
$r2 = $d.Content
$r2.Find.Execute(" as a high first. Achieved official first class results in second and third year.")
$r2.Text = " as a high first. Achieved official first-class results in second year and third year."

$r2b = $d.Content
$r2b.Find.Execute(" as a high first. Achieved official first-class results in second year and third year.")
$start2 = $r2b.Start

$seg2 = $d.Range($start2 + 41, $start2 + 42)   # "-"
$seg2.Font.Bold = $true
$seg2.Font.Bold = $false

$seg3 = $d.Range($start2 + 42, $start2 + 66)   # "class results in second "
$seg3.Font.Bold = $true
$seg3.Font.Bold = $false

$seg4 = $d.Range($start2 + 66, $start2 + 71)   # "year "
$seg4.Font.Bold = $true
$seg4.Font.Bold = $false

$seg5 = $d.Range($start2 + 71, $start2 + 85)   # "and third year"
$seg5.Font.Bold = $true
$seg5.Font.Bold = $false

$seg6 = $d.Range($start2 + 85, $start2 + 86)   # "."
$seg6.Font.Bold = $true
$seg6.Font.Bold = $false

# =====================================================================
# Change 3: " (2019 before COVID)" (bold)
#   -> " (2019" | "," | " before COVID)"
#   (re-isolates the untouched "4 A* or equivalent" / "." neighbour runs
#    that the engine would otherwise coalesce into this edit)
# =====================================================================
$r3 = $d.Content
$r3.Find.Execute(" (2019 before COVID)")
$r3.Text = " (2019, before COVID)"

$r3b = $d.Content
$r3b.Find.Execute("4 A* or equivalent (2019, before COVID).")
$start3 = $r3b.Start

$c1 = $d.Range($start3 + 18, $start3 + 24)     # " (2019"
$c1.Font.Bold = $false
$c1.Font.Bold = $true

$c2 = $d.Range($start3 + 24, $start3 + 25)     # ","
$c2.Font.Bold = $false
$c2.Font.Bold = $true

$c3 = $d.Range($start3 + 25, $start3 + 39)     # " before COVID)"
$c3.Font.Bold = $false
$c3.Font.Bold = $true

$c4 = $d.Range($start3 + 39, $start3 + 40)     # "."
$c4.Font.Bold = $false
$c4.Font.Bold = $true

Write-Output "Done"
